$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '98.785.80'
$ws.Range("E2").Value = '  +1.48%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.484.50'
$ws.Range("E3").Value = '  +5.97%  '

# Row 4
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '262.24'
$ws.Range("E5").Value = '  +3.49%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '675.69'
$ws.Range("E6").Value = '  +9.10%  '

# Row 7
$ws.Range("E7").Value = '  +8.80%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.465'
$ws.Range("E8").Value = '  +17.82%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.11'
$ws.Range("E9").Value = '  +23.40%  '

# Row 10
$ws.Range("E10").Value = '  -0.06%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.479.13'
$ws.Range("E11").Value = '  +5.92%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.221'
$ws.Range("E12").Value = '  +11.67%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.82'
$ws.Range("E13").Value = '  +12.32%  '

# Row 14
$ws.Range("E14").Value = '  +12.61%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.30'
$ws.Range("E15").Value = '  +15.97%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.137.90'
$ws.Range("E16").Value = '  +5.97%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '98.980.73'
$ws.Range("E17").Value = '  +1.82%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.22'
$ws.Range("E18").Value = '  +34.36%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.482.06'
$ws.Range("E19").Value = '  +5.80%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.61'
$ws.Range("E20").Value = '  +17.43%  '

# Row 21
$ws.Range("E21").Value = '  +4.66%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '536.46'
$ws.Range("E22").Value = '  +12.93%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.80'
$ws.Range("E23").Value = '  +15.95%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000220'
$ws.Range("E24").Value = '  +9.48%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.454'
$ws.Range("E25").Value = '  +52.21%  '

# Row 26
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '102.70'
$ws.Range("E26").Value = '  +17.06%  '

# Row 27
$ws.Range("B27").Value = 'NEARProtocol'
$ws.Range("C27").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.41'
$ws.Range("E27").Value = '  +15.42%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.00'
$ws.Range("E28").Value = '  +10.96%  '

# Row 29
$ws.Range("E29").Value = '  +15.31%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.202'
$ws.Range("E30").Value = '  +9.52%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.45'
$ws.Range("E31").Value = '  +17.24%  '

# Row 32
$ws.Range("E32").Value = '  +0.15%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.591'
$ws.Range("E33").Value = '  +31.30%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.95'
$ws.Range("E34").Value = '  +13.32%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.41%  '

# Row 36
$ws.Range("E36").Value = '  +16.43%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.03'
$ws.Range("E37").Value = '  +12.94%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.163'
$ws.Range("E38").Value = '  +11.23%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '542.32'
$ws.Range("E39").Value = '  +11.58%  '

# Row 40
$ws.Range("E40").Value = '  +16.23%  '

# Row 41
$ws.Range("E41").Value = '  +0.00%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.871'
$ws.Range("E42").Value = '  +10.30%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0438'
$ws.Range("E43").Value = '  +36.15%  '

# Row 44
$ws.Range("E44").Value = '  +13.01%  '

# Row 45
$ws.Range("E45").Value = '  +1.22%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.24'
$ws.Range("E46").Value = '  +17.49%  '

# Row 47
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.13'
$ws.Range("E47").Value = '  +13.62%  '

# Row 48
$ws.Range("B48").Value = 'ImmutableX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.61'
$ws.Range("E48").Value = '  +20.35%  '

# Row 49
$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  +0.02%  '

# Row 50
$ws.Range("E50").Value = '  +16.01%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.79'
$ws.Range("E51").Value = '  +14.13%  '
